# "April 5 update 2" -- add two new computed columns to the right of the
# existing data: NewNonRentCost (AI) and NewRentCost (AJ).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New headers in row 1 (style copied from AH1 so they match the other
# bold/centered header cells).
$ws.Range("AI1").Value = "NewNonRentCost"
$ws.Range("AJ1").Value = "NewRentCost"
$ws.Range("AH1").Copy() | Out-Null
$ws.Range("AI1:AJ1").PasteSpecial(-4122) | Out-Null

# Row 2 holds the "master" formulas; rows 3:25 are filled down from them,
# which is what turns them into one shared-formula group per column.
$ws.Range("AI2").Formula = "=AH2+AG2"
$ws.Range("AJ2").Formula = "=AH2+AG2+AE2"

$ws.Range("AI3:AI25").Formula = "=AH3+AG3"
$ws.Range("AJ3:AJ25").Formula = "=AH3+AG3+AE3"

# Leave the selection on AJ2, as in the saved workbook.
$ws.Range("AJ2").Select() | Out-Null
